$d = $word.ActiveDocument

# The four new bullet points to append, in order. The last one keeps the
# trailing "_GoBack" bookmark that currently sits right after the final
# run of the document.
$newTexts = @(
    "Assigning to the master the new tile from the queue was placed was wrong. Whenever we take a new tile from the queue, we refresh the table. It makes sense because it is actually a new search. ",
    "Then, check for the parent tiles of the newly taken tile from the queue and place them on the board. Once placed, place the new tile. Then only check if we need to work with the new tile. ",
    "Removed the stack of stationary files since it didn’t really make sense. We only push and pop one element. So, I changed it to a TileTree object.",
    "ERROR: The tiles were not being placed correctly. Once the available location was found, lets say (5,4), the column never started from 0 again even for next row. So, fixed it by changing the value of starting point of column to 0."
)

# Drop the existing "_GoBack" bookmark - it needs to move to the end of the
# new last paragraph, and it's easier to re-create it than to relocate it.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Build the text to insert: the bullets separated by paragraph marks, plus
# one throw-away trailing character. Inserting right at the end of the
# document and then placing a bookmark exactly at content-end is unreliable
# in this host, so we park a sentinel character after everything, bookmark
# right before it (a perfectly ordinary mid-document position), then trim
# the sentinel back off.
$joined = [string]::Join("`r", $newTexts)
$sentinel = "@"

$end = $d.Content.End
$insertRange = $d.Range($end, $end)
$insertRange.InsertAfter("`r" + $joined + $sentinel)

$lastPara = $d.Paragraphs.Last
$bookmarkPos = $lastPara.Range.End - 2
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinelRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$sentinelRange.Delete()
